# Refresh the crypto price (column D) and 1h volume-change (column E)
# figures on Sheet1, rows 2-51, per the GitHub Actions data refresh run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.207.73"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "'1.671.76"
$ws.Range("E3").Value = "  -1.56%  "

$ws.Range("E4").Value = "  -0.65%  "

$ws.Range("D5").Value = "'211.60"
$ws.Range("E5").Value = "  -3.18%  "

$ws.Range("D6").Value = "'0.5281"
$ws.Range("E6").Value = "  -3.50%  "

$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("D8").Value = "'0.2641"
$ws.Range("E8").Value = "  -3.53%  "

$ws.Range("D9").Value = "'0.06286"
$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("D10").Value = "'21.32"
$ws.Range("E10").Value = "  -2.99%  "

$ws.Range("D11").Value = "'0.07555"
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("D12").Value = "'1.669.35"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("D13").Value = "'4.447"
$ws.Range("E13").Value = "  -2.43%  "

$ws.Range("D14").Value = "'0.5602"
$ws.Range("E14").Value = "  -4.35%  "

$ws.Range("D15").Value = "'67.13"
$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("D16").Value = "'0.000008003"
$ws.Range("E16").Value = "  -4.79%  "

$ws.Range("D17").Value = "'26.249.29"
$ws.Range("E17").Value = "  -0.89%  "

$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "'4.788"
$ws.Range("E19").Value = "  -3.21%  "

$ws.Range("D20").Value = "'187.44"
$ws.Range("E20").Value = "  -2.18%  "

$ws.Range("E21").Value = "  -5.25%  "

$ws.Range("E22").Value = "  -1.17%  "

$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").Value = "'149.60"
$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("D25").Value = "'0.1259"
$ws.Range("E25").Value = "  -3.87%  "

$ws.Range("D26").Value = "'7.578"
$ws.Range("E26").Value = "  -4.48%  "

$ws.Range("D27").Value = "'15.95"
$ws.Range("E27").Value = "  +0.92%  "

$ws.Range("D28").Value = "'0.06179"
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").Value = "'1.367"
$ws.Range("E29").Value = "  -1.82%  "

$ws.Range("D30").Value = "'1.283"
$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("D31").Value = "'3.500"
$ws.Range("E31").Value = "  -3.14%  "

$ws.Range("D32").Value = "'3.428"
$ws.Range("E32").Value = "  -4.67%  "

$ws.Range("E33").Value = "  -3.39%  "

$ws.Range("D34").Value = "'1.000"
$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("D35").Value = "'0.6074"
$ws.Range("E35").Value = "  -1.45%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "'2.737"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").Value = "'6.125"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("E39").Value = "  -2.36%  "

$ws.Range("D40").Value = "'1.098.34"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").Value = "'0.8786"
$ws.Range("E41").Value = "  -0.48%  "

$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("D44").Value = "'1.822.42"
$ws.Range("E44").Value = "  -1.47%  "

$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = "  +2.03%  "

$ws.Range("D46").Value = "'55.90"
$ws.Range("E46").Value = "  -3.04%  "

$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("D48").Value = "'8.010"
$ws.Range("E48").Value = "  -2.70%  "

$ws.Range("D49").Value = "'0.05229"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").Value = "'5.980"
$ws.Range("E51").Value = "  -2.41%  "
